$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set all data values in B2:D9 to 0 (area check / specs files update)
$ws.Range("B2:D9").Value = 0
